$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107; this shifts the existing rows 107..227
# down to 108..228 (matching the diff's net effect: dimension A1:R227 -> A1:R228).
$ws.Rows(107).Insert()

# Populate the newly inserted row 107 with a fresh weekly data record.
# The "constant" columns (A,B,C,E,F,G,H,I,R) are identical on every data
# row in this sheet, so we replicate them here too.
$ws.Range("A107").Value = 10
$ws.Range("B107").Value = "Vega Modelo de Temuco"
$ws.Range("C107").Value = "La Araucanía"
$ws.Range("D107").Value = 44679
$ws.Range("E107").Value = 9
$ws.Range("F107").Value = 100112052
$ws.Range("G107").Value = "Albahaca"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 50
$ws.Range("K107").Value = 5000
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = 5000
$ws.Range("N107").Value = "$/paquete"
$ws.Range("O107").Value = "Región de Arica y Parinacota"
$ws.Range("P107").Value = 5000
$ws.Range("Q107").Value = 1
$ws.Range("R107").Value = "Hortaliza"
